$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptos list: updated Price (col D) and Volume(1h) (col E) figures,
# plus two rank swaps (Algorand/TheSandbox at rows 42-43, Decentraland/Quant at rows 46-47).
# Price cells that look like plain numbers get NumberFormat "@" applied first so Excel
# keeps them as literal text (preserving trailing zeros) instead of coercing to a number.
$ws.Range("D2").Value = '27.443.26'
$ws.Range("E2").Value = '  -1.09%  '
$ws.Range("D3").Value = '1.833.29'
$ws.Range("E3").Value = '  -0.93%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.004'
$ws.Range("E4").Value = '  -2.64%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '315.22'
$ws.Range("E5").Value = '  -2.11%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.004'
$ws.Range("E6").Value = '  -2.51%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4297'
$ws.Range("E7").Value = '  -2.18%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3704'
$ws.Range("E8").Value = '  -2.60%  '
$ws.Range("E9").Value = '  -1.78%  '
$ws.Range("E10").Value = '  -1.97%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '21.18'
$ws.Range("E11").Value = '  -1.69%  '
$ws.Range("D12").Value = '1.838.77'
$ws.Range("E12").Value = '  -0.76%  '
$ws.Range("E13").Value = '  +0.12%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.361'
$ws.Range("E14").Value = '  -2.62%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.07091'
$ws.Range("E15").Value = '  -1.09%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '88.13'
$ws.Range("E16").Value = '  +3.55%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.007'
$ws.Range("E17").Value = '  -2.91%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008919'
$ws.Range("E18").Value = '  -1.73%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.004'
$ws.Range("E19").Value = '  -2.49%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.24'
$ws.Range("E20").Value = '  -1.61%  '
$ws.Range("D21").Value = '27.452.50'
$ws.Range("E21").Value = '  -1.12%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.170'
$ws.Range("E22").Value = '  -1.99%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.92'
$ws.Range("E23").Value = '  -2.97%  '
$ws.Range("D24").Value = '2.064.49'
$ws.Range("E24").Value = '  -1.06%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.001'
$ws.Range("E25").Value = '  -2.44%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '153.62'
$ws.Range("E26").Value = '  -3.04%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.50'
$ws.Range("E27").Value = '  -1.06%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.147'
$ws.Range("E28").Value = '  +7.65%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.282'
$ws.Range("E29").Value = '  -1.10%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '117.30'
$ws.Range("E30").Value = '  -0.43%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08859'
$ws.Range("E31").Value = '  -2.32%  '
$ws.Range("E32").Value = '  -0.25%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7687'
$ws.Range("E33").Value = '  -0.48%  '
$ws.Range("E34").Value = '  -1.88%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.911'
$ws.Range("E35").Value = '  -2.94%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.004'
$ws.Range("E36").Value = '  -2.60%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.120'
$ws.Range("E37").Value = '  -2.70%  '
$ws.Range("E38").Value = '  -0.68%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05294'
$ws.Range("E39").Value = '  +0.19%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.187'
$ws.Range("E40").Value = '  +4.63%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.870'
$ws.Range("E41").Value = '  +0.46%  '
$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5087'
$ws.Range("E42").Value = '  -1.83%  '
$ws.Range("B43").Value = 'Algorand'
$ws.Range("C43").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1676'
$ws.Range("E43").Value = '  +0.23%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.679'
$ws.Range("E44").Value = '  -0.69%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.64'
$ws.Range("E45").Value = '  -0.59%  '
$ws.Range("B46").Value = 'Quant'
$ws.Range("C46").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '106.39'
$ws.Range("E46").Value = '  -3.67%  '
$ws.Range("B47").Value = 'Decentraland'
$ws.Range("C47").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4742'
$ws.Range("E47").Value = '  +0.76%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.06426'
$ws.Range("E48").Value = '  -2.21%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.004'
$ws.Range("E49").Value = '  -2.82%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.668'
$ws.Range("E50").Value = '  -2.24%  '
$ws.Range("E51").Value = '  -3.26%  '
